# Roboflow Annotation Report 7/29/2025
# Append the next day's progress row ("29/7/2031") to the tracking table
# on Sheet1, mirroring the values already logged for the previous day
# (28/7/2030) on row 67, then move the table/window to reflect where the
# author was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Table1 currently spans D4:J67 (63 data rows). Ask the ListObject to grow
# by one row - this keeps the table's ref/autoFilter and the sheet
# dimension in sync automatically.
$lo = $ws.ListObjects.Item("Table1")
$lo.ListRows.Add() | Out-Null

# New row index is one past the old last data row.
$newRowIndex = $lo.Range.Row + $lo.Range.Rows.Count - 1

# Clone the formatting of the previous data row (row 67) onto the new row
# (68) before writing values, so fonts/borders/number formats match the
# rest of the table exactly.
$prevRowIndex = $newRowIndex - 1
$ws.Range("D" + $prevRowIndex + ":J" + $prevRowIndex).Copy()
$ws.Range("D" + $newRowIndex + ":J" + $newRowIndex).PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Rows.Item($newRowIndex).RowHeight = $ws.Rows.Item($prevRowIndex).RowHeight

# Fill in the new day's numbers (same counts as the prior entry, just a
# fresh date stamp).
$ws.Range("D" + $newRowIndex).Value = "29/7/2031"
$ws.Range("E" + $newRowIndex).Value = 380
$ws.Range("F" + $newRowIndex).Value = 950
$ws.Range("G" + $newRowIndex).Value = 0
$ws.Range("H" + $newRowIndex).Value = 0
$ws.Range("I" + $newRowIndex).Value = 1012
$ws.Range("J" + $newRowIndex).Value = "N/A"

# Scroll the sheet down to the new rows and leave the selection where the
# author ended up.
$win = $excel.ActiveWindow
$win.ScrollRow = 45
$win.ScrollColumn = 2
$ws.Range("G78").Select()
